# Updated backlog with more probable dates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sprint End Date (row 4) pushed out
$ws.Range("G4").Value = "2/4/2026"

# Network Architecture Review sub-tasks (rows 16-19)
$ws.Range("G16").Value = "1/30/2026"
$ws.Range("G17").Value = "2/1/2026"
$ws.Range("G19").Value = "2/2/2026"

# Vulnerability Identification Report sub-tasks (rows 20-23)
$ws.Range("G20").Value = "1/30/2026"
$ws.Range("G21").Value = "2/1/2026"
$ws.Range("G23").Value = "2/2/2026"

# Best Practice Comparison sub-tasks (rows 24-27)
$ws.Range("G24").Value = "1/30/2026"
$ws.Range("G25").Value = "2/1/2026"
$ws.Range("G27").Value = "2/2/2026"

# Threat Scenario Analysis sub-tasks (rows 28-31)
$ws.Range("G28").Value = "1/30/2026"
$ws.Range("G29").Value = "2/1/2026"
$ws.Range("G31").Value = "2/2/2026"

# Improvement Recommendations sub-tasks (rows 32-35)
$ws.Range("G32").Value = "1/30/2026"
$ws.Range("G33").Value = "2/1/2026"
$ws.Range("G35").Value = "2/2/2026"

# Executive Summary sub-tasks (rows 36-39)
$ws.Range("G36").Value = "2/2/2026"
$ws.Range("G37").Value = "2/3/2026"
$ws.Range("G38").ClearContents()
$ws.Range("E39").Value = "Shubham"
$ws.Range("G39").Value = "2/4/2026"

# Slide Deck sub-tasks (rows 40-41)
$ws.Range("G40").Value = "2/4/2026"

# User Stories sub-tasks (rows 42-43)
$ws.Range("G42").Value = "2/2/2026"
$ws.Range("G43").Value = "2/2/2026"

# Update the active selection to match the authored session state
$ws.Activate()
$ws.Range("J41").Select()
